$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 (I0) and J1 (IF) - copy formatting from the existing H1
# header (bold/bordered/centered style) and then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-39.
$i0 = @(8,9,5,8,7,8,5,6,4,4,6,8,5,6,5,6,6,7,6,7,8,9,9,8,7,6,7,9,8,6,8,6,9,8,7,6,2,5)
$iF = @(8,9,6,8,7,9,5,6,4,5,6,8,6,7,5,6,6,7,6,7,8,9,9,8,7,6,7,9,8,6,8,7,9,9,7,6,3,5)

for ($r = 0; $r -lt $i0.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $iF[$r]
}
